$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 5 (extra data rows), keeping only header row and first data row
$ws.Range("A3:F5").EntireRow.Delete()

# Update row 2 with the new values
$ws.Range("A2").Value = "Multijugador"
$ws.Range("B2").Value = "Ulises"
$ws.Range("C2").Value = 395
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = "Derrota"
$ws.Range("F2").Value = "2025-11-24 21:51:16"
